$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country order: Finlandia now sorts before Gabon ---
# Row 102 used to be Gabon, now becomes Finlandia (with refreshed stats)
# Row 103 used to be Finlandia, now becomes Gabon (stats unchanged from old Gabon row)
$ws.Range("A102").Value = "Finlandia"
$ws.Range("A103").Value = "Gabon"

# --- Swap country order: Eslovaquia now sorts before Malaui ---
# Row 110 used to be Malaui, now becomes Eslovaquia (with refreshed stats)
# Row 111 used to be Eslovaquia, now becomes Malaui (stats unchanged from old Malaui row)
$ws.Range("A110").Value = "Eslovaquia"
$ws.Range("A111").Value = "Malaui"

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 11:38"

# --- Refresh numeric statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 24: Filipinas
$ws.Range("B24").Value = 269407
$ws.Range("C24").Value = 3544
$ws.Range("D24").Value = 207352
$ws.Range("E24").Value = 57392
$ws.Range("G24").Value = 34
$ws.Range("H24").Value = 4663

# Row 26: Indonesia
$ws.Range("B26").Value = 225030
$ws.Range("C26").Value = 3507
$ws.Range("D26").Value = 161065
$ws.Range("E26").Value = 55000
$ws.Range("G26").Value = 124
$ws.Range("H26").Value = 8965

# Row 27: Israel
$ws.Range("B27").Value = 162273
$ws.Range("C27").Value = 1905
$ws.Range("D27").Value = 120443
$ws.Range("E27").Value = 40689
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 1141

# Row 48: Polonia
$ws.Range("B48").Value = 75134
$ws.Range("C48").Value = 605
$ws.Range("D48").Value = 61548
$ws.Range("E48").Value = 11359
$ws.Range("G48").Value = 24
$ws.Range("H48").Value = 2227

# Row 70: Austria
$ws.Range("B70").Value = 34305
$ws.Range("C70").Value = 764
$ws.Range("D70").Value = 27354
$ws.Range("E70").Value = 6194

# Row 97: Malasia
$ws.Range("B97").Value = 9969
$ws.Range("C97").Value = 23
$ws.Range("D97").Value = 9209
$ws.Range("E97").Value = 632

# Row 102: now Finlandia (refreshed stats)
$ws.Range("B102").Value = 8725
$ws.Range("C102").Value = 98
$ws.Range("D102").Value = 7500
$ws.Range("E102").Value = 888
$ws.Range("H102").Value = 337

# Row 103: now Gabon (keeps the old Gabon stats)
$ws.Range("B103").Value = 8654
$ws.Range("D103").Value = 7785
$ws.Range("E103").Value = 816
$ws.Range("H103").Value = 53

# Row 110: now Eslovaquia (refreshed stats)
$ws.Range("B110").Value = 5768
$ws.Range("C110").Value = 188
$ws.Range("D110").Value = 3214
$ws.Range("E110").Value = 2516
$ws.Range("H110").Value = 38

# Row 111: now Malaui (keeps the old Malaui stats)
$ws.Range("B111").Value = 5697
$ws.Range("D111").Value = 3742
$ws.Range("E111").Value = 1777
$ws.Range("H111").Value = 178

# Row 126: Eslovenia
$ws.Range("B126").Value = 3831
$ws.Range("C126").Value = 82
$ws.Range("D126").Value = 2789
$ws.Range("E126").Value = 907
